{"js": "// Locate the paragraphs that hold the \"<<judgeRecital>>\" and\n// \"<<hearingOrder>>\" merge fields and replace them:\n//  - the judgeRecital paragraph becomes the new recital sentence\n//  - the hearingOrder paragraph (and one of the blank paragraphs that\n//    used to separate the two) is removed entirely\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet judgeIdx = -1;\nlet hearingIdx = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (judgeIdx === -1 && t.indexOf(\"<<judgeRecital>>\") !== -1) {\n    judgeIdx = i;\n  }\n  if (hearingIdx === -1 && t.indexOf(\"<<hearingOrder>>\") !== -1) {\n    hearingIdx = i;\n  }\n}\n\nif (judgeIdx === -1 || hearingIdx === -1) {\n  throw new Error(\"Could not locate judgeRecital/hearingOrder paragraphs\");\n}\n\nconst judgePara = paragraphs.items[judgeIdx];\nconst hearingPara = paragraphs.items[hearingIdx];\n\n// Replace the whole text of the judgeRecital paragraph with the new\n// recital sentence (merge fields stay as literal << >> placeholders).\njudgePara.insertText(\n  \"Upon the application of <<applicantName>> dated <<applicationDate>> and upon considering the information provided by the parties:\",\n  \"Replace\"\n);\n\n// The hearingOrder paragraph together with the blank paragraph right\n// after it collapses down to a single blank paragraph, so delete the\n// hearingOrder paragraph outright.\nhearingPara.delete();\n\nawait context.sync();\n", "ps1": "# Replace the \"<<judgeRecital>>\" merge-field paragraph with the new\n# recital sentence, then remove the \"<<hearingOrder>>\" merge-field\n# paragraph entirely (one of the blank separator paragraphs collapses\n# away along with it).\n\n$d = $word.ActiveDocument\n\n$newRecital = \"Upon the application of <<applicantName>> dated <<applicationDate>> and upon considering the information provided by the parties:\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"<<judgeRecital>>\"\n$find.Replacement.Text = $newRecital\n$find.Execute(\n    \"<<judgeRecital>>\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    $newRecital,\n    2\n)\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -match \"<<hearingOrder>>\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
